# Delete the two blank spacer rows (12 and 13) that separate the "Week 1"
# and "Week 2" sprint-planning tables. Deleting the entire rows shifts the
# "Week 2" table (previously rows 15-28) up by two rows (to rows 13-26),
# matching the rest of the diff exactly (the cell contents/styles of that
# table are otherwise untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("12:13").Delete() | Out-Null

# Reflect the scrolled/selected state recorded in the saved file.
$ws.Range("D20").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 12
